$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 4, shifting existing rows 4-18 down to 5-19.
$ws.Rows("4:4").Insert()

# Populate the new row 4 with the latest weekly price entry.
$ws.Range("A4").Value = 10
$ws.Range("B4").Value = "Vega Modelo de Temuco"
$ws.Range("C4").Value = "La Araucanía"
$ws.Range("D4").Value = 44707
$ws.Range("E4").Value = 9
$ws.Range("F4").Value = 100112042
$ws.Range("G4").Value = "Locoto"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 4700
$ws.Range("L4").Value = 4700
$ws.Range("M4").Value = 4700
$ws.Range("N4").Value = "`$/kilo"
$ws.Range("O4").Value = "Región de Arica y Parinacota"
$ws.Range("P4").Value = 4700
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = "Hortaliza"
